$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "1. The system must support Aadhaar/PAN card OCR and verification in compliance with UIDAI and Income Tax Department regulations."
$ws.Range("A4").Value = "2. The system must ensure secure integration with credit bureaus (CIBIL, Experian) as per the Credit Information Companies (Regulation) Act, 2005."
$ws.Range("A5").Value = "3. The AI-based eligibility and risk assessment module must comply with the principles of fairness, transparency, and non-discrimination as per the AI Ethics Guidelines."
$ws.Range("A6").Value = "4. Real-time KYC document verification (OCR + face match) must be in compliance with the Reserve Bank of India (RBI) KYC guidelines."
$ws.Range("A7").Value = "5. The system must ensure secure OTP-based authentication for loan sanction and e-agreement signing in line with the Information Technology Act, 2000."
$ws.Range("A8").Value = "6. The system must ensure secure integration with core banking for disbursement as per the Banking Regulation Act, 1949."
$ws.Range("A9").Value = "7. The E-sign integration using DigiLocker + Aadhaar eKYC must comply with the provisions of the Information Technology Act, 2000 and the Aadhaar Act, 2016."
$ws.Range("A10").Value = "8. The system must ensure data privacy and protection in line with the Personal Data Protection Bill, 2019."
$ws.Range("A11").Value = "9. The system must support English and Hindi languages as per the Official Languages Act, 1963."
$ws.Range("A12").Value = "10. The system must be prepared to adapt to compliance changes due to upcoming RBI guidelines."
$ws.Range("A13").Value = "11. The system must ensure that loan disbursement is only through verified bank accounts as per the Prevention of Money Laundering Act, 2002."
$ws.Range("A14").Value = "12. The system must be designed to be accessible to all users, including those with disabilities, in compliance with the Rights of Persons with Disabilities Act, 2016."
